# Rename the "_old"/"_new" header-name suffixes to the respective
# format-version suffixes ("_FV2404"/"_FV2410"), then turn the used
# range into a native Excel Table (ListObject) with a frozen header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters A..U in header order, paired with their new names.
$headerNames = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404",
    "diff",
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt $headerNames.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headerNames[$i]
}

# Convert A1:U67 into a native table (adds xl/tables/table1.xml,
# the tableParts wiring and the auto filter).
$usedRange = $ws.Range("A1:U67")
$tbl = $ws.ListObjects.Add(1, $usedRange, $null, 1)
$tbl.Name = "Table1"

# Freeze the header row (pane split after row 1).
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
